$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the stray note in D44 ("change clusters in db to eigendocs + OPTICS") entirely
$ws.Range("D44").Clear() | Out-Null

# 2. Insert the two new log rows (50 & 51) that follow the existing 27.09.2023 entry (row 49).
#    Copy formatting (borders/number-format) from row 49's A:B cells first, then fill in values.
$ws.Range("A49:B49").Copy() | Out-Null
$ws.Range("A50:B50").PasteSpecial(-4122) | Out-Null
$ws.Range("A51:B51").PasteSpecial(-4122) | Out-Null

$ws.Range("A50").Value = 45197
$ws.Range("B50").Value2 = "BA: TFIDF, Doc2Vec"

$ws.Range("A51").Value = 45198
$ws.Range("B51").Value2 = ""

$ws.Rows.Item(50).RowHeight = 18
$ws.Rows.Item(51).RowHeight = 18

# 3. Restore the view state (scroll/selection) to where the author left off.
$ws.Range("D55").Select() | Out-Null

Write-Output "edit applied"
